$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = [double]"1.273760039130352E-07"
$ws.Range("E2").Value = [double]"1.273760039130352E-07"

# Row 3
$ws.Range("D3").Value = [double]"0.001536591649188994"
$ws.Range("E3").Value = [double]"0.001536591649188994"

# Row 4
$ws.Range("D4").Value = [double]"6.264361322180282E-05"
$ws.Range("E4").Value = [double]"6.264361322180282E-05"

# Row 5
$ws.Range("D5").Value = [double]"1.321622470031253E-05"
$ws.Range("E5").Value = [double]"1.321622470031253E-05"

# Row 6
$ws.Range("D6").Value = [double]"0.1870592056908366"
$ws.Range("E6").Value = [double]"0.1870592056908366"

# Row 7
$ws.Range("D7").Value = [double]"0.8356824989492203"
$ws.Range("E7").Value = [double]"0.1643175010507797"

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = [double]"7.312781674127107E-07"
$ws.Range("E8").Value = [double]"0.9999992687218325"

# Row 9
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = [double]"5.341450409448754E-16"
$ws.Range("E9").Value = [double]"0.9999999999999994"

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = [double]"1.469325468107369E-06"
$ws.Range("E10").Value = [double]"0.9999985306745319"

# Row 11
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = [double]"0.2578318332462012"
$ws.Range("E11").Value = [double]"0.7421681667537988"
$ws.Range("F11").Value = [double]"6.44687032699585"
$ws.Range("G11").Value = [double]"0.6"

# Row 12
$ws.Range("D12").Value = [double]"1.484727933813741E-10"
$ws.Range("E12").Value = [double]"1.484727933813741E-10"

# Row 13
$ws.Range("D13").Value = [double]"1.991111457275314E-07"
$ws.Range("E13").Value = [double]"1.991111457275314E-07"

# Row 14
$ws.Range("D14").Value = [double]"1.603137247098961E-06"
$ws.Range("E14").Value = [double]"1.603137247098961E-06"

# Row 15
$ws.Range("D15").Value = [double]"7.403026299060435E-08"
$ws.Range("E15").Value = [double]"7.403026299060435E-08"

# Row 16
$ws.Range("D16").Value = [double]"0.05183561844560355"
$ws.Range("E16").Value = [double]"0.05183561844560355"

# Row 17
$ws.Range("D17").Value = [double]"0.8757795608878421"
$ws.Range("E17").Value = [double]"0.1242204391121579"

# Row 18
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = [double]"1.64497084213964E-11"
$ws.Range("E18").Value = [double]"0.9999999999835503"

# Row 19
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = [double]"1.827945652468018E-24"
$ws.Range("E19").Value = 1

# Row 20
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = [double]"1.656877932017172E-08"
$ws.Range("E20").Value = [double]"0.9999999834312207"

# Row 21
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = [double]"0.01512657650101828"
$ws.Range("E21").Value = [double]"0.9848734234989818"
$ws.Range("F21").Value = [double]"10.17824745178223"
$ws.Range("G21").Value = [double]"0.6"
